# Update Leave Card 1/19/2024 4:44 pm
# -----------------------------------------------------------------------
# This script reproduces, via Excel COM interop, the edit that was made
# to "COLETO, HANY ROY.xlsx" Sheet1:
#   * A new leave entry (VL 23 days, 4/1-5/1/2024) is recorded in the row
#     that used to be the blank "April" placeholder row (row 184).
#   * A new blank "May" placeholder row is inserted right after it,
#     partially filled in with an SP entry (B185) and a literal end date
#     in K185, which pushes the rest of the yearly placeholder rows
#     (185-213, with their column-A month dates) down by one row.
#   * The Table1 range, and the worksheet, both grow by one row (the
#     former "last row" of the table - with its special bottom border
#     styling - moves from row 243 to the new row 244).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# -------------------------------------------------------------------
# 1) Grow the table by one row first, so that writing into row 244
#    later on is recognised as "inside the table" (keeps structured
#    references such as Table1[[#This Row],[EARNED]] working).
# -------------------------------------------------------------------
$tbl.Resize($ws.Range("A8:K244"))

# -------------------------------------------------------------------
# 2) Re-create the old "last row" of the table (row 243) as the new
#    last row (row 244), carrying over its special bottom-border
#    formatting and its EARNED-helper formula.
# -------------------------------------------------------------------
$ws.Range("A243:K243").Copy()
$ws.Range("A244:K244").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("G244").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Row 243 becomes a normal interior row (copy formatting from row 242).
$ws.Range("A242:K242").Copy()
$ws.Range("A243:K243").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# -------------------------------------------------------------------
# 3) Shift the column-A placeholder month dates for rows 185-213 down
#    into rows 186-214, and blank out A185 (it becomes the row that is
#    being turned into the new leave entry).
# -------------------------------------------------------------------
$datesToShift = $ws.Range("A185:A213").Value2
$ws.Range("A186:A214").Value2 = $datesToShift
$ws.Range("A185").ClearContents()

# -------------------------------------------------------------------
# 4) Fill in the new leave-card entries.
# -------------------------------------------------------------------
# Row 184: VL(23-0-0) covering 23 days, 4/1 - 5/1/2024
$ws.Range("B184").Value2 = "VL(23-0-0)"
$ws.Range("D184").Value2 = 23
$ws.Range("K184").Value2 = "4/1 - 5/1/2024"

# Row 185: SP(1-0-0) entry; K185 carries a literal end date (45323)
# instead of the usual formula, matching the other manually dated
# entries elsewhere in the sheet (e.g. K87, K121, K171).
$ws.Range("B185").Value2 = "SP(1-0-0)"
$ws.Range("G171").Copy()
$ws.Range("G185").PasteSpecial(-4122)        # xlPasteFormats (style 42 -> 13)
$excel.CutCopyMode = 0
$ws.Range("K171").Copy()
$ws.Range("K185").PasteSpecial(-4122)        # xlPasteFormats (style 20 -> 49)
$excel.CutCopyMode = 0
$ws.Range("K185").Value2 = 45323

# -------------------------------------------------------------------
# 5) Update the worksheet selection to match where the user ended up.
# -------------------------------------------------------------------
$ws.Range("K185").Select()
